# Automated map update (2025-09-22 07:32:13)
# Inserts a new record as row 71 on sheet "NEW", pushing the existing
# rows 71 ("2485" / LA PLATA AV. 1095) and 72 ("232" / Gorostiaga 2286)
# down to rows 72 and 73 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 71 (and everything below it) down by one row.
$ws.Rows("71:71").Insert()

# --- New row 71 ------------------------------------------------------
# Columns A, B and D look numeric/date-like ("2711", "9/22/2025", "12")
# but must be stored as plain text, matching every other row in this
# sheet. Force text via NumberFormat, assign, then clear the format so
# no new cell style is left behind (keeps the cell on the default
# style, like its neighbours).
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "2711"
$ws.Range("A71").ClearFormats()

$ws.Range("B71").NumberFormat = "@"
$ws.Range("B71").Value = "9/22/2025"
$ws.Range("B71").ClearFormats()

$ws.Range("C71").Value = "RUIZ HUIDOBRO 3620"

$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "12"
$ws.Range("D71").ClearFormats()

$ws.Range("E71").Value = "Pendiente ADM"
$ws.Range("F71").Value = "NEW"
$ws.Range("G71").Value = "Pendiente"
$ws.Range("H71").Value = "Trapaso de redes y desmonte"
$ws.Range("I71").Value = 1
$ws.Range("J71").Value = "Desmonte"
$ws.Range("K71").Value = "Sin equipos"
$ws.Range("L71").Value = "Pasante"
$ws.Range("M71").Value = -58.484082
$ws.Range("N71").Value = -34.549702
$ws.Range("O71").Value = "Saavedra"
$ws.Range("P71").Value = "Capital Norte"
